# ---------------------------------------------------------------------------
# Scheduled Sheets refresh: re-pulls current Market Board averages/prices for
# each Leve row and recomputes the NQ/HQ profit columns (H:N) per job sheet.
# Generated by the scheduled runner; values below mirror the refreshed feed.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# === Sheet: ALC ===
$ws = $wb.Worksheets.Item(1)
# Row 40
$ws.Range("H40").Value = 40001980
$ws.Range("I40").Value = 1425
$ws.Range("J40").Value = 47621132
$ws.Range("K40").Value = 1425
$ws.Range("L40").Value = 47621132
$ws.Range("M40").Value = -1250
$ws.Range("N40").Value = -47621482

# Row 51
$ws.Range("H51").Value = 2584.6667
$ws.Range("I51").Value = 2585
$ws.Range("J51").Value = 2584.2856
$ws.Range("K51").Value = 2585
$ws.Range("L51").Value = 2584.2856
$ws.Range("M51").Value = -2101
$ws.Range("N51").Value = -3552.2856

# Row 53
$ws.Range("H53").Value = 308.7
$ws.Range("I53").Value = 118
$ws.Range("K53").Value = 118
$ws.Range("M53").Value = 519

# Row 64
$ws.Range("H64").Value = 27450
$ws.Range("I64").Value = 35333.332
$ws.Range("J64").Value = 3800
$ws.Range("K64").Value = 35333.332
$ws.Range("L64").Value = 3800
$ws.Range("M64").Value = -35085.332
$ws.Range("N64").Value = -4296

# Row 67
$ws.Range("H67").Value = 27450
$ws.Range("I67").Value = 35333.332
$ws.Range("J67").Value = 3800
$ws.Range("K67").Value = 35333.332
$ws.Range("L67").Value = 3800
$ws.Range("M67").Value = -34475.332
$ws.Range("N67").Value = -5516

# Row 74
$ws.Range("H74").Value = 12981656
$ws.Range("I74").Value = 12981656
$ws.Range("K74").Value = 12981656
$ws.Range("M74").Value = -12980720

# Row 77
$ws.Range("H77").Value = 12981656
$ws.Range("I77").Value = 12981656
$ws.Range("K77").Value = 64908280
$ws.Range("M77").Value = -64903600

# Row 106
$ws.Range("H106").Value = 2296.5
$ws.Range("I106").Value = 2134.1538
$ws.Range("K106").Value = 2134.1538
$ws.Range("M106").Value = -1503.1538

# Row 121
$ws.Range("H121").Value = 815.7273
$ws.Range("J121").Value = 817.3
$ws.Range("L121").Value = 2451.9
$ws.Range("N121").Value = -5945.9

# Row 129
$ws.Range("H129").Value = 1544258.4
$ws.Range("J129").Value = 1544258.4
$ws.Range("L129").Value = 4632775.199999999
$ws.Range("N129").Value = -4642775.199999999

# Row 137
$ws.Range("H137").Value = 6758420.5
$ws.Range("I137").Value = 1308.8959
$ws.Range("J137").Value = 19233088
$ws.Range("K137").Value = 3926.6877
$ws.Range("L137").Value = 57699264
$ws.Range("M137").Value = -1376.6877
$ws.Range("N137").Value = -57704364

# === Sheet: ARM ===
$ws = $wb.Worksheets.Item(2)
# Row 32
$ws.Range("H32").Value = 21111.037
$ws.Range("I32").Value = 18591.783
$ws.Range("J32").Value = 58396
$ws.Range("K32").Value = 18591.783
$ws.Range("L32").Value = 58396
$ws.Range("M32").Value = -18304.783
$ws.Range("N32").Value = -58970

# Row 132
$ws.Range("H132").Value = 2067
$ws.Range("I132").Value = 1467.9333
$ws.Range("K132").Value = 4403.7999
$ws.Range("M132").Value = -1873.7999

# === Sheet: BSM ===
$ws = $wb.Worksheets.Item(3)
# Row 68
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# Row 69
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

# Row 71
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# Row 72
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

# Row 82
$ws.Range("H82").Value = 15752.333

# Row 85
$ws.Range("H85").Value = 15752.333

# Row 94
$ws.Range("H94").Value = 3064.5
$ws.Range("I94").Value = 2849.4443
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 2849.4443
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -2398.4443
$ws.Range("N94").Value = -5902

# Row 106
$ws.Range("H106").Value = 29000
$ws.Range("J106").Value = 29000
$ws.Range("L106").Value = 29000
$ws.Range("N106").Value = -31524

# Row 107
$ws.Range("H107").Value = 45016.668
$ws.Range("I107").Value = 1125
$ws.Range("J107").Value = 132800
$ws.Range("K107").Value = 1125
$ws.Range("L107").Value = 132800
$ws.Range("M107").Value = 795
$ws.Range("N107").Value = -136640

# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# Row 109
$ws.Range("H109").Value = 34800
$ws.Range("J109").Value = 34800
$ws.Range("L109").Value = 34800
$ws.Range("N109").Value = -37574

# Row 134
$ws.Range("H134").Value = 52022.773
$ws.Range("I134").Value = 2097.2188
$ws.Range("J134").Value = 251725
$ws.Range("K134").Value = 6291.6564
$ws.Range("L134").Value = 755175
$ws.Range("M134").Value = -3756.6564
$ws.Range("N134").Value = -760245

# === Sheet: CRP ===
$ws = $wb.Worksheets.Item(4)
# Row 62
$ws.Range("H62").Value = 9250
$ws.Range("I62").Value = 8000
$ws.Range("J62").Value = 9875
$ws.Range("K62").Value = 8000
$ws.Range("L62").Value = 9875
$ws.Range("M62").Value = -7376
$ws.Range("N62").Value = -11123

# Row 65
$ws.Range("H65").Value = 9250
$ws.Range("I65").Value = 8000
$ws.Range("J65").Value = 9875
$ws.Range("K65").Value = 40000
$ws.Range("L65").Value = 49375
$ws.Range("M65").Value = -36880
$ws.Range("N65").Value = -55615

# Row 107
$ws.Range("H107").Value = 409.04166
$ws.Range("J107").Value = 394.3
$ws.Range("L107").Value = 394.3
$ws.Range("N107").Value = -4234.3

# Row 132
$ws.Range("H132").Value = 3537.2727
$ws.Range("I132").Value = 1002.4
$ws.Range("K132").Value = 3007.2
$ws.Range("M132").Value = -477.1999999999998

# Row 140
$ws.Range("H140").Value = 50331.11
$ws.Range("J140").Value = 50331.11
$ws.Range("L140").Value = 50331.11
$ws.Range("N140").Value = -60691.11

# === Sheet: CUL ===
$ws = $wb.Worksheets.Item(5)
# Row 68
$ws.Range("H68").Value = 1266.2222
$ws.Range("I68").Value = 653.6842
$ws.Range("J68").Value = 1647.8032
$ws.Range("K68").Value = 1961.0526
$ws.Range("L68").Value = 4943.4096
$ws.Range("M68").Value = -1150.0526
$ws.Range("N68").Value = -6565.4096

# Row 71
$ws.Range("H71").Value = 1266.2222
$ws.Range("I71").Value = 653.6842
$ws.Range("J71").Value = 1647.8032
$ws.Range("K71").Value = 5883.1578
$ws.Range("L71").Value = 14830.2288
$ws.Range("M71").Value = -1827.1578
$ws.Range("N71").Value = -22942.2288

# Row 113
$ws.Range("H113").Value = 512.0123
$ws.Range("I113").Value = 346.42856
$ws.Range("J113").Value = 546.61194
$ws.Range("K113").Value = 1039.28568
$ws.Range("L113").Value = 1639.83582
$ws.Range("M113").Value = 1130.71432
$ws.Range("N113").Value = -5979.83582

# Row 134
$ws.Range("H134").Value = 1725.2354
$ws.Range("I134").Value = 1725.2354
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5175.706200000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -105.7062000000005
$ws.Range("N134").ClearContents()

# === Sheet: GSM ===
$ws = $wb.Worksheets.Item(6)
# Row 70
$ws.Range("H70").Value = 4356.3477
$ws.Range("I70").Value = 4084.0527
$ws.Range("J70").Value = 5649.75
$ws.Range("K70").Value = 4084.0527
$ws.Range("L70").Value = 5649.75
$ws.Range("M70").Value = -3814.0527
$ws.Range("N70").Value = -6189.75

# Row 73
$ws.Range("H73").Value = 4356.3477
$ws.Range("I73").Value = 4084.0527
$ws.Range("J73").Value = 5649.75
$ws.Range("K73").Value = 4084.0527
$ws.Range("L73").Value = 5649.75
$ws.Range("M73").Value = -3148.0527
$ws.Range("N73").Value = -7521.75

# Row 107
$ws.Range("H107").Value = 624.8
$ws.Range("I107").Value = 536
$ws.Range("J107").Value = 980
$ws.Range("K107").Value = 536
$ws.Range("L107").Value = 980
$ws.Range("M107").Value = 1384
$ws.Range("N107").Value = -4820

# Row 138
$ws.Range("H138").Value = 36419
$ws.Range("J138").Value = 36419
$ws.Range("L138").Value = 36419
$ws.Range("N138").Value = -46699

# === Sheet: LTW ===
$ws = $wb.Worksheets.Item(7)
# Row 46
$ws.Range("H46").Value = 472.72726
$ws.Range("I46").Value = 433.33334
$ws.Range("K46").Value = 433.33334
$ws.Range("M46").Value = -245.33334

# Row 68
$ws.Range("H68").Value = 2340
$ws.Range("I68").Value = 2333.3333
$ws.Range("J68").Value = 2350
$ws.Range("K68").Value = 2333.3333
$ws.Range("L68").Value = 2350
$ws.Range("M68").Value = -1584.3333
$ws.Range("N68").Value = -3848

# Row 71
$ws.Range("H71").Value = 2340
$ws.Range("I71").Value = 2333.3333
$ws.Range("J71").Value = 2350
$ws.Range("K71").Value = 11666.6665
$ws.Range("L71").Value = 11750
$ws.Range("M71").Value = -7922.666499999999
$ws.Range("N71").Value = -19238

# Row 132
$ws.Range("H132").Value = 2528112.8
$ws.Range("I132").Value = 3137446.8
$ws.Range("J132").Value = 3728.2856
$ws.Range("K132").Value = 9412340.399999999
$ws.Range("L132").Value = 11184.8568
$ws.Range("M132").Value = -9409810.399999999
$ws.Range("N132").Value = -16244.8568

# Row 139
$ws.Range("H139").Value = 47715
$ws.Range("J139").Value = 47715
$ws.Range("L139").Value = 47715
$ws.Range("N139").Value = -57995
